# Update betting odds values in the "Jogos da Semana" worksheet
# to reflect the latest FlashScore snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 3.2
$ws.Range("I2").Value = 4.75
$ws.Range("K2").Value = 1.95
$ws.Range("O2").Value = 1.5
$ws.Range("P2").Value = 2.5
$ws.Range("Q2").Value = 2.6
$ws.Range("R2").Value = 1.48
$ws.Range("S2").Value = 1.57
$ws.Range("T2").Value = 2.25
$ws.Range("U2").Value = 2.25
$ws.Range("V2").Value = 1.57
$ws.Range("X2").Value = 7.5
$ws.Range("AA2").Value = 19
$ws.Range("AC2").Value = 6.5
$ws.Range("AE2").Value = 21
$ws.Range("AF2").Value = 81
$ws.Range("AG2").Value = 9.5
$ws.Range("AK2").Value = 41
$ws.Range("AR2").Value = 67
$ws.Range("AT2").Value = 2.25
$ws.Range("AU2").Value = 9.5
$ws.Range("AW2").Value = 6
$ws.Range("AZ2").Value = 101
$ws.Range("BA2").Value = 151
$ws.Range("BB2").Value = 451

# Row 3
$ws.Range("AS3").Value = 151
$ws.Range("AZ3").Value = 126
$ws.Range("BB3").Value = 351

# Row 5
$ws.Range("Q5").Value = 2.08
$ws.Range("R5").Value = 1.73

# Row 6
$ws.Range("G6").Value = 2.15
$ws.Range("I6").Value = 3.3
$ws.Range("J6").Value = 2.88
$ws.Range("K6").Value = 2.05
$ws.Range("O6").Value = 1.36
$ws.Range("P6").Value = 3
$ws.Range("U6").Value = 1.91
$ws.Range("V6").Value = 1.8
$ws.Range("Z6").Value = 19
$ws.Range("AA6").Value = 19
$ws.Range("AC6").Value = 8.5
$ws.Range("AF6").Value = 51
$ws.Range("AG6").Value = 9
$ws.Range("AO6").Value = 12
$ws.Range("AP6").Value = 23
$ws.Range("AY6").Value = 29

# Row 7
$ws.Range("G7").Value = 4.33
$ws.Range("I7").Value = 1.85
$ws.Range("J7").Value = 5.5
$ws.Range("U7").Value = 2.25
$ws.Range("V7").Value = 1.57
$ws.Range("AS7").Value = 451
$ws.Range("AW7").Value = 3.6

# Row 10
$ws.Range("K10").Value = 1.95

# Row 11
$ws.Range("G11").Value = 2.05
$ws.Range("I11").Value = 3.9
$ws.Range("J11").Value = 2.88
$ws.Range("L11").Value = 4.5
$ws.Range("M11").Value = 1.1
$ws.Range("N11").Value = 7
$ws.Range("X11").Value = 8.5
$ws.Range("AD11").Value = 6
$ws.Range("AL11").Value = 41
$ws.Range("AN11").Value = 4
$ws.Range("AO11").Value = 12
$ws.Range("AW11").Value = 5.5
$ws.Range("AY11").Value = 34
